$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (copying row 19's formatting down), which
# pushes the old "Total" row from 21 down to 22.
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(20).Insert()

# Fill in the new entry that records the "Updated prices" task
$ws.Range("A20").Value = "Updated prices"
$ws.Range("B20").Value = 43364
$ws.Range("C20").Value = 1

# The Total row moved from 21 down to 22 because of the inserted row; leave row 21 blank
# and update the Total row's sum formula to include the newly inserted row.
$ws.Range("C22").Formula = "=SUM(C2:C21)"

$ws.Range("C22").Select()
